$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.200.34"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.603.57"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'212.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.484"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.249"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "'0.0612"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "'18.20"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").Value = "'0.0814"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "1.828.74"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.604.30"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "26.183.38"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "'61.78"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D20").Value = "'200.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'5.99"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +2.18%  "
$ws.Range("D25").Value = "'144.11"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("D28").Value = "'15.18"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "'6.55"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("E30").Value = "  +3.80%  "
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").Value = "'2.92"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("D36").Value = "1.163.01"
$ws.Range("E36").Value = "  +5.05%  "
$ws.Range("E37").Value = "  +3.70%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "'0.784"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "'0.780"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").Value = "'5.30"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.94%  "
$ws.Range("D44").Value = "1.740.34"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "'91.51"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").Value = "'53.98"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "'0.408"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "0.0₇0944"
$ws.Range("E51").Value = "  +5.60%  "
